# Updates the API link in the read-me:
#   - "hi" becomes the intro line "API's for testing:"
#   - a blank line, a hyperlink to the REST API docs, and two trailing
#     blank lines are appended
#   - a "Hyperlink" character style (blue + underline) backs the new link,
#     the same way Word auto-creates it the first time a hyperlink is added

$d = $word.ActiveDocument

# Create the "Hyperlink" character style up front so the run we add below
# can reference it via rStyle, and so it round-trips into styles.xml.
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = $d.Styles("DefaultParagraphFont")
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Underline = 1
$hlStyle.Font.Color = 16711680

# "hi" -> "API's for testing:" (straight quote auto-corrects to the
# typographic U+2019 apostrophe, matching real Word behavior).
$d.Content.Find.Execute("hi", $true, $false, $false, $false, $false, $true, `
                         1, $false, "API's for testing:", 2)

# Reserve the trailing paragraphs: blank, link-holder, blank, blank.
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertXML("<w:p $wns/><w:p $wns/><w:p $wns/><w:p $wns/>")

# Turn the 3rd (newly added) paragraph into the hyperlink paragraph.
$linkRange = $d.Paragraphs(3).Range
$linkRange.Collapse(0)
$d.Hyperlinks.Add($linkRange, "http://rest.learncode.academy/", $null, $null, `
                   "http://rest.learncode.academy/") | Out-Null

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
